$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -14.15280000000001
$ws.Range("D4").Value = -8.572400000000002
$ws.Range("D5").Value = -8.5322
$ws.Range("C6").Value = -12.0301
$ws.Range("C7").Value = -12.1055
$ws.Range("D8").Value = -8.0496
$ws.Range("C16").Value = -11.7162
$ws.Range("D16").Value = -8.436900000000009
$ws.Range("C20").Value = -14.7293
$ws.Range("D22").Value = -7.8094
